$wb = $excel.ActiveWorkbook

# Sheet ALC, row 17 (diff hunk 0)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1026.9207
$ws.Range("J17").Value = 865.0484
$ws.Range("L17").Value = 2595.1452
$ws.Range("N17").Value = -2931.1452

# Sheet ALC, row 86 (diff hunk 1)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 1773.3334
$ws.Range("I86").Value = 1700
$ws.Range("K86").Value = 1700
$ws.Range("M86").Value = -577

# Sheet ALC, row 89 (diff hunk 2)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 1773.3334
$ws.Range("I89").Value = 1700
$ws.Range("K89").Value = 8500
$ws.Range("M89").Value = -2884

# Sheet ALC, row 98 (diff hunk 3)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H98").Value = 2130.2173
$ws.Range("I98").Value = 2337.5
$ws.Range("J98").Value = 748.3333
$ws.Range("K98").Value = 2337.5
$ws.Range("L98").Value = 748.3333
$ws.Range("M98").Value = -839.5
$ws.Range("N98").Value = -3744.3333

# Sheet ALC, row 122 (diff hunk 4)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H122").Value = 2130.2173
$ws.Range("I122").Value = 2337.5
$ws.Range("J122").Value = 748.3333
$ws.Range("K122").Value = 7012.5
$ws.Range("L122").Value = 2244.9999
$ws.Range("M122").Value = -4562.5
$ws.Range("N122").Value = -7144.9999

# Sheet ALC, row 132 (diff hunk 5)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H132").Value = 1054.7778
$ws.Range("I132").Value = 969.7646999999999
$ws.Range("K132").Value = 2909.2941
$ws.Range("M132").Value = -379.2941000000001

# Sheet ALC, row 137 (diff hunk 6)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 1563.3572
$ws.Range("I137").Value = 1468.7
$ws.Range("J137").Value = 1800
$ws.Range("K137").Value = 4406.1
$ws.Range("L137").Value = 5400
$ws.Range("M137").Value = -1856.1
$ws.Range("N137").Value = -10500

# Sheet ALC, row 138 (diff hunk 7)
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 3763.6775
$ws.Range("J138").Value = 3646.9333
$ws.Range("L138").Value = 10940.7999
$ws.Range("N138").Value = -21220.7999

# Sheet ARM, row 32 (diff hunk 8)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 2913.6086
$ws.Range("I32").Value = 2084.7415
$ws.Range("K32").Value = 2084.7415
$ws.Range("M32").Value = -1797.7415

# Sheet ARM, row 61 (diff hunk 9)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 4481.7617
$ws.Range("I61").Value = 2841.2
$ws.Range("K61").Value = 2841.2
$ws.Range("M61").Value = -2629.2

# Sheet ARM, row 74 (diff hunk 10)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H74").Value = 1285.6904
$ws.Range("I74").Value = 847.4211
$ws.Range("J74").Value = 5449.25
$ws.Range("K74").Value = 847.4211
$ws.Range("L74").Value = 5449.25
$ws.Range("M74").Value = 26.57889999999998
$ws.Range("N74").Value = -7197.25

# Sheet ARM, row 77 (diff hunk 11)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H77").Value = 1285.6904
$ws.Range("I77").Value = 847.4211
$ws.Range("J77").Value = 5449.25
$ws.Range("K77").Value = 4237.1055
$ws.Range("L77").Value = 27246.25
$ws.Range("M77").Value = 130.8945000000003
$ws.Range("N77").Value = -35982.25

# Sheet ARM, row 102 (diff hunk 12)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 1894.3334
$ws.Range("I102").Value = 1647.3077
$ws.Range("K102").Value = 1647.3077
$ws.Range("M102").Value = -25.30770000000007

# Sheet ARM, row 136 (diff hunk 13)
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 4481.7617
$ws.Range("I136").Value = 2841.2
$ws.Range("K136").Value = 8523.599999999999
$ws.Range("M136").Value = -5973.599999999999

# Sheet BSM, row 20 (diff hunk 14)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H20").Value = 1780.1052
$ws.Range("I20").Value = 1636
$ws.Range("J20").Value = 2320.5
$ws.Range("K20").Value = 1636
$ws.Range("L20").Value = 2320.5
$ws.Range("M20").Value = -1389
$ws.Range("N20").Value = -2814.5

# Sheet BSM, row 76 (diff hunk 15)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H76").Value = 55156.5
$ws.Range("J76").Value = 55156.5
$ws.Range("L76").Value = 55156.5
$ws.Range("N76").Value = -55786.5

# Sheet BSM, row 79 (diff hunk 16)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H79").Value = 55156.5
$ws.Range("J79").Value = 55156.5
$ws.Range("L79").Value = 55156.5
$ws.Range("N79").Value = -57340.5

# Sheet BSM, row 86 (diff hunk 17)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 102304.1
$ws.Range("I86").Value = 2184.5
$ws.Range("K86").Value = 2184.5
$ws.Range("M86").Value = -1061.5

# Sheet BSM, row 89 (diff hunk 18)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H89").Value = 102304.1
$ws.Range("I89").Value = 2184.5
$ws.Range("K89").Value = 10922.5
$ws.Range("M89").Value = -5306.5

# Sheet BSM, row 107 (diff hunk 19)
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 3251.8333
$ws.Range("I107").Value = 3251.8333
$ws.Range("K107").Value = 3251.8333
$ws.Range("M107").Value = -1331.8333

# Sheet CRP, row 62 (diff hunk 20)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H62").Value = 2477.5
$ws.Range("J62").Value = 2975
$ws.Range("L62").Value = 2975
$ws.Range("N62").Value = -4223

# Sheet CRP, row 65 (diff hunk 21)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H65").Value = 2477.5
$ws.Range("J65").Value = 2975
$ws.Range("L65").Value = 14875
$ws.Range("N65").Value = -21115

# Sheet CRP, row 70 (diff hunk 22)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H70").Value = 20000
$ws.Range("J70").Value = 20000
$ws.Range("L70").Value = 20000
$ws.Range("N70").Value = -20630

# Sheet CRP, row 73 (diff hunk 23)
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H73").Value = 20000
$ws.Range("J73").Value = 20000
$ws.Range("L73").Value = 20000
$ws.Range("N73").Value = -22184

# Sheet CUL, row 108 (diff hunk 24)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H108").Value = 2002.75
$ws.Range("I108").Value = 2002.75
$ws.Range("K108").Value = 6008.25
$ws.Range("M108").Value = -3128.25

# Sheet CUL, row 131 (diff hunk 25)
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H131").Value = 749.99
$ws.Range("J131").Value = 781.94446
$ws.Range("L131").Value = 2345.83338
$ws.Range("N131").Value = -12425.83338

# Sheet GSM, row 97 (diff hunk 26)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 845.4167
$ws.Range("I97").Value = 782.4375
$ws.Range("K97").Value = 782.4375
$ws.Range("M97").Value = -286.4375

# Sheet GSM, row 122 (diff hunk 27)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H122").Value = 1993.1538
$ws.Range("I122").Value = 1584
$ws.Range("J122").Value = 2343.8572
$ws.Range("K122").Value = 4752
$ws.Range("L122").Value = 7031.571599999999
$ws.Range("M122").Value = -2302
$ws.Range("N122").Value = -11931.5716

# Sheet GSM, row 126 (diff hunk 28)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 34828.676
$ws.Range("I126").Value = 2775.0908
$ws.Range("K126").Value = 8325.2724
$ws.Range("M126").Value = -5855.2724

# Sheet GSM, row 132 (diff hunk 29)
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H132").Value = 3181.2173
$ws.Range("I132").Value = 2788.45
$ws.Range("J132").Value = 5799.6665
$ws.Range("K132").Value = 8365.349999999999
$ws.Range("L132").Value = 17398.9995
$ws.Range("M132").Value = -5835.349999999999
$ws.Range("N132").Value = -22458.9995

# Sheet LTW, row 70 (diff hunk 30)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H70").Value = 30000
$ws.Range("J70").Value = 30000
$ws.Range("L70").Value = 30000
$ws.Range("N70").Value = -30540

# Sheet LTW, row 73 (diff hunk 31)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H73").Value = 30000
$ws.Range("J73").Value = 30000
$ws.Range("L73").Value = 30000
$ws.Range("N73").Value = -31872

# Sheet LTW, row 75 (diff hunk 32)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H75").Value = 30000
$ws.Range("J75").Value = 30000
$ws.Range("L75").Value = 30000
$ws.Range("N75").Value = -31872

# Sheet LTW, row 78 (diff hunk 33)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H78").Value = 30000
$ws.Range("J78").Value = 30000
$ws.Range("L78").Value = 90000
$ws.Range("N78").Value = -99360

# Sheet LTW, row 132 (diff hunk 34)
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 2174
$ws.Range("I132").Value = 2014.1111
$ws.Range("K132").Value = 6042.3333
$ws.Range("M132").Value = -3512.3333

# Sheet WVR, row 14 (diff hunk 35)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H14").Value = 3840.5
$ws.Range("I14").Value = 10000
$ws.Range("J14").Value = 2608.6
$ws.Range("K14").Value = 10000
$ws.Range("L14").Value = 2608.6
$ws.Range("M14").Value = -9832
$ws.Range("N14").Value = -2944.6

# Sheet WVR, row 64 (diff hunk 36)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H64").Value = 29999
$ws.Range("J64").Value = 29999
$ws.Range("L64").Value = 29999
$ws.Range("N64").Value = -30495

# Sheet WVR, row 67 (diff hunk 37)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H67").Value = 29999
$ws.Range("J67").Value = 29999
$ws.Range("L67").Value = 29999
$ws.Range("N67").Value = -31715

# Sheet WVR, row 70 (diff hunk 38)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H70").Value = 29039.4
$ws.Range("J70").Value = 29039.4
$ws.Range("L70").Value = 29039.4
$ws.Range("N70").Value = -29669.4

# Sheet WVR, row 73 (diff hunk 39)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H73").Value = 29039.4
$ws.Range("J73").Value = 29039.4
$ws.Range("L73").Value = 29039.4
$ws.Range("N73").Value = -31223.4

# Sheet WVR, row 82 (diff hunk 40)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H82").Value = 61767
$ws.Range("J82").Value = 61767
$ws.Range("L82").Value = 61767
$ws.Range("N82").Value = -62533

# Sheet WVR, row 85 (diff hunk 41)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H85").Value = 61767
$ws.Range("J85").Value = 61767
$ws.Range("L85").Value = 61767
$ws.Range("N85").Value = -64419

# Sheet WVR, row 122 (diff hunk 42)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H122").Value = 71934.55
$ws.Range("I122").Value = 98121.875
$ws.Range("K122").Value = 294365.625
$ws.Range("M122").Value = -291915.625

# Sheet WVR, row 132 (diff hunk 43)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H132").Value = 1899.5555
$ws.Range("I132").Value = 867.6667
$ws.Range("J132").Value = 2415.5
$ws.Range("K132").Value = 2603.0001
$ws.Range("L132").Value = 7246.5
$ws.Range("M132").Value = -73.0001000000002
$ws.Range("N132").Value = -12306.5

# Sheet WVR, row 136 (diff hunk 44)
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H136").Value = 2716.1667
$ws.Range("I136").Value = 2878.8125
$ws.Range("J136").Value = 2530.2856
$ws.Range("K136").Value = 8636.4375
$ws.Range("L136").Value = 7590.8568
$ws.Range("M136").Value = -6086.4375
$ws.Range("N136").Value = -12690.8568
